# Update the cryptos price/volume table (columns D and E, rows 2-51).
# Values are stored as text in the sheet (prices use "." as thousands
# separators so they are not valid numbers, e.g. "69.007.38"); a leading
# apostrophe is used for the few D-column values that Excel would otherwise
# auto-convert to a number (stripping meaningful trailing/leading zeros),
# forcing them to stay text exactly as in the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.007.38"
$ws.Range("E2").Value = "  -3.37%  "
$ws.Range("D3").Value = "3.494.39"
$ws.Range("E3").Value = "  -5.62%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'576.89"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "'171.05"
$ws.Range("E6").Value = "  -3.97%  "
$ws.Range("D7").Value = "3.486.27"
$ws.Range("E7").Value = "  -5.65%  "
$ws.Range("D8").Value = "'0.604"
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -5.98%  "
$ws.Range("D11").Value = "'6.47"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "'0.579"
$ws.Range("E12").Value = "  -4.81%  "
$ws.Range("D13").Value = "'46.40"
$ws.Range("E13").Value = "  -5.63%  "
$ws.Range("D14").Value = "'0.0000272"
$ws.Range("E14").Value = "  -5.05%  "
$ws.Range("D15").Value = "4.056.94"
$ws.Range("E15").Value = "  -5.81%  "
$ws.Range("D16").Value = "'8.49"
$ws.Range("E16").Value = "  -5.57%  "
$ws.Range("D17").Value = "'614.09"
$ws.Range("E17").Value = "  -9.38%  "
$ws.Range("D18").Value = "68.954.21"
$ws.Range("E18").Value = "  -3.65%  "
$ws.Range("D19").Value = "3.484.00"
$ws.Range("E19").Value = "  -6.29%  "
$ws.Range("D20").Value = "'0.122"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "'17.24"
$ws.Range("E21").Value = "  -3.96%  "
$ws.Range("D22").Value = "'11.07"
$ws.Range("E22").Value = "  -4.37%  "
$ws.Range("D23").Value = "'0.880"
$ws.Range("E23").Value = "  -6.64%  "
$ws.Range("D24").Value = "'15.80"
$ws.Range("E24").Value = "  -9.49%  "
$ws.Range("D25").Value = "'96.76"
$ws.Range("E25").Value = "  -5.62%  "
$ws.Range("E26").Value = "  -5.42%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "'2.62"
$ws.Range("E28").Value = "  -7.20%  "
$ws.Range("D29").Value = "'9.27"
$ws.Range("E29").Value = "  -10.43%  "
$ws.Range("D30").Value = "'32.48"
$ws.Range("E30").Value = "  -7.69%  "
$ws.Range("D31").Value = "'3.14"
$ws.Range("E31").Value = "  -8.36%  "
$ws.Range("D32").Value = "'8.45"
$ws.Range("E32").Value = "  -7.52%  "
$ws.Range("E33").Value = "  -8.92%  "
$ws.Range("D34").Value = "'6.93"
$ws.Range("E34").Value = "  -5.75%  "
$ws.Range("D35").Value = "'626.87"
$ws.Range("E35").Value = "  +6.06%  "
$ws.Range("D36").Value = "'10.68"
$ws.Range("E36").Value = "  -4.61%  "
$ws.Range("E37").Value = "  -5.76%  "
$ws.Range("D38").Value = "'3.42"
$ws.Range("E38").Value = "  -15.84%  "
$ws.Range("D39").Value = "'56.46"
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'0.0443"
$ws.Range("E41").Value = "  -3.22%  "
$ws.Range("E42").Value = "  -6.28%  "
$ws.Range("D43").Value = "3.347.51"
$ws.Range("E43").Value = "  -8.92%  "
$ws.Range("D44").Value = "'0.325"
$ws.Range("E44").Value = "  -6.48%  "
$ws.Range("D45").Value = "'32.58"
$ws.Range("E45").Value = "  -7.67%  "
$ws.Range("D46").Value = "0.0₃0690"
$ws.Range("E46").Value = "  -10.20%  "
$ws.Range("D47").Value = "'2.56"
$ws.Range("E47").Value = "  -7.85%  "
$ws.Range("D48").Value = "'2.76"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("D49").Value = "'0.129"
$ws.Range("E49").Value = "  -3.09%  "
$ws.Range("D50").Value = "'132.00"
$ws.Range("E50").Value = "  -2.97%  "
$ws.Range("E51").Value = "  +13.67%  "
